$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2
Set-TextValue "D2" "29.379.69"
$ws.Range("E2").Value = "  +0.05%  "

# Row 3
Set-TextValue "D3" "1.848.04"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4
Set-TextValue "D4" "0.9999"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
Set-TextValue "D5" "240.21"
$ws.Range("E5").Value = "  +0.13%  "

# Row 6
Set-TextValue "D6" "0.6296"
$ws.Range("E6").Value = "  -0.58%  "

# Row 7
Set-TextValue "D7" "1.001"
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
Set-TextValue "D8" "0.07605"
$ws.Range("E8").Value = "  +0.81%  "

# Row 9
Set-TextValue "D9" "0.2927"
$ws.Range("E9").Value = "  -1.03%  "

# Row 10
Set-TextValue "D10" "24.42"
$ws.Range("E10").Value = "  -1.19%  "

# Row 11
Set-TextValue "D11" "0.07740"
$ws.Range("E11").Value = "  +0.12%  "

# Row 12
Set-TextValue "D12" "1.846.39"
$ws.Range("E12").Value = "  -6.98%  "

# Row 13
Set-TextValue "D13" "0.00001100"
$ws.Range("E13").Value = "  +11.28%  "

# Row 14
Set-TextValue "D14" "5.000"

# Row 15
Set-TextValue "D15" "0.6782"
$ws.Range("E15").Value = "  -0.69%  "

# Row 16
Set-TextValue "D16" "83.57"
$ws.Range("E16").Value = "  +0.58%  "

# Row 17
Set-TextValue "D17" "2.105.83"
$ws.Range("E17").Value = "  -7.00%  "

# Row 18
$ws.Range("E18").Value = "  +0.09%  "

# Row 19
Set-TextValue "D19" "29.398.20"
$ws.Range("E19").Value = "  +0.00%  "

# Row 20
Set-TextValue "D20" "228.37"
$ws.Range("E20").Value = "  -0.88%  "

# Row 21
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
Set-TextValue "D22" "1.001"
$ws.Range("E22").Value = "  +0.11%  "

# Row 23
Set-TextValue "D23" "7.480"
$ws.Range("E23").Value = "  -0.81%  "

# Row 24
$ws.Range("E24").Value = "  +0.13%  "

# Row 25
$ws.Range("E25").Value = "  +0.69%  "

# Row 26
Set-TextValue "D26" "0.1396"
$ws.Range("E26").Value = "  -0.60%  "

# Row 27
Set-TextValue "D27" "8.341"
$ws.Range("E27").Value = "  -0.24%  "

# Row 28
$ws.Range("E28").Value = "  -0.23%  "

# Row 29
Set-TextValue "D29" "1.462"
$ws.Range("E29").Value = "  -0.27%  "

# Row 30
Set-TextValue "D30" "1.300"
$ws.Range("E30").Value = "  +3.87%  "

# Row 31
Set-TextValue "D31" "0.05591"
$ws.Range("E31").Value = "  -2.01%  "

# Row 32
Set-TextValue "D32" "4.103"
$ws.Range("E32").Value = "  -0.40%  "

# Row 33
Set-TextValue "D33" "4.028"
$ws.Range("E33").Value = "  +0.04%  "

# Row 34
Set-TextValue "D34" "1.842"
$ws.Range("E34").Value = "  -0.02%  "

# Row 35
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
Set-TextValue "D36" "0.7108"
$ws.Range("E36").Value = "  -0.58%  "

# Row 37
Set-TextValue "D37" "2.586"
$ws.Range("E37").Value = "  -0.23%  "

# Row 38
Set-TextValue "D38" "1.238.92"
$ws.Range("E38").Value = "  -0.86%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.01801"
$ws.Range("E39").Value = "  -0.48%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D40" "2.770"
$ws.Range("E40").Value = "  -1.04%  "

# Row 41
Set-TextValue "D41" "6.410"
$ws.Range("E41").Value = "  +5.20%  "

# Row 42
Set-TextValue "D42" "0.9049"
$ws.Range("E42").Value = "  +0.39%  "

# Row 43
$ws.Range("E43").Value = "  +0.07%  "

# Row 44
Set-TextValue "D44" "101.69"
$ws.Range("E44").Value = "  -0.12%  "

# Row 45
$ws.Range("E45").Value = "  -0.63%  "

# Row 46
$ws.Range("E46").Value = "  +2.14%  "

# Row 47
Set-TextValue "D47" "7.143"
$ws.Range("E47").Value = "  +1.05%  "

# Row 48
Set-TextValue "D48" "0.4019"
$ws.Range("E48").Value = "  +0.12%  "

# Row 49
Set-TextValue "D49" "8.993"
$ws.Range("E49").Value = "  -1.32%  "

# Row 50
Set-TextValue "D50" "1.677"
$ws.Range("E50").Value = "  -1.41%  "

# Row 51
Set-TextValue "D51" "0.1120"
$ws.Range("E51").Value = "  -0.32%  "
